$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add "Groupes CM" column -------------------------------------------
# A new column is inserted before the old column F ("Heures TD"), so the
# former F:I ("Heures TD", "Groupes TD", "Heures TP", "Groupes TP") shift
# right into G:J, and the new column F holds the "Groupes CM" figures.
$ws.Columns("F:F").Insert()

$ws.Range("F1").Value2 = "Groupes CM"

# "Groupes CM" value per course row (number of CM groups, mirroring how
# "Groupes TD"/"Groupes TP" record the TD/TP group counts).
$groupesCM = @{
  2  = 1
  3  = 1
  4  = 3
  5  = 1
  6  = 1
  7  = 1
  8  = 1
  9  = 1
  10 = 1
  11 = 1
  12 = 3
  13 = 1
  14 = 1
  15 = 1
}

foreach ($r in $groupesCM.Keys) {
  $ws.Cells.Item($r, 6).Value2 = $groupesCM[$r]
}

# --- Update view state ---------------------------------------------------
# Scroll the window so column B becomes the left-most visible column, and
# move the active selection to D15.
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$ws.Range("D15").Select()
